$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Update the Date property
$wsMeta.Range("B8").Value = "2026-01-01T13:37:23+00:00"

# Update the Description property (both on Metadata sheet and Elements sheet)
$newDescription = "Extension to link goal evaluation observations to the patient goals being evaluated."
$wsMeta.Range("B11").Value = $newDescription
$wsElem.Range("M2").Value = $newDescription

# Update the Extension.value[x] Reference text (rename onc-patient-goal -> onc-nursing-goal)
$wsElem.Range("K6").Value = "Reference(https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/onc-nursing-goal)`n"

# The column got a hair wider after the re-generation (bestFit recompute); nudge it
# back towards the recorded width as closely as this engine's column-width rounding allows.
$wsElem.Columns.Item(11).ColumnWidth = 77.5
